$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '66.365.73'
Set-TextValue 'E2' '  -6.10%  '
Set-TextValue 'D3' '3.195.51'
Set-TextValue 'E3' '  -9.41%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '573.41'
Set-TextValue 'E5' '  -6.68%  '
Set-TextValue 'D6' '148.56'
Set-TextValue 'E6' '  -14.59%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  +0.07%  '
Set-TextValue 'D8' '3.185.47'
Set-TextValue 'E8' '  -9.52%  '
Set-TextValue 'D9' '0.538'
Set-TextValue 'E9' '  -11.74%  '
Set-TextValue 'E10' '  -14.60%  '
Set-TextValue 'E11' '  -12.15%  '
Set-TextValue 'D12' '0.493'
Set-TextValue 'E12' '  -16.32%  '
Set-TextValue 'D13' '38.01'
Set-TextValue 'E13' '  -18.39%  '
Set-TextValue 'D14' '0.0000239'
Set-TextValue 'E14' '  -13.74%  '
Set-TextValue 'D15' '3.701.25'
Set-TextValue 'E15' '  -9.68%  '
Set-TextValue 'D16' '66.389.98'
Set-TextValue 'E16' '  -6.06%  '
Set-TextValue 'D17' '3.191.87'
Set-TextValue 'E17' '  -9.58%  '
Set-TextValue 'E18' '  -6.63%  '
Set-TextValue 'D19' '529.88'
Set-TextValue 'E19' '  -14.06%  '
Set-TextValue 'D20' '7.06'
Set-TextValue 'E20' '  -16.60%  '
Set-TextValue 'D21' '14.90'
Set-TextValue 'E21' '  -16.23%  '
Set-TextValue 'D22' '0.750'
Set-TextValue 'E22' '  -15.30%  '
Set-TextValue 'D23' '7.64'
Set-TextValue 'E23' '  -15.12%  '
Set-TextValue 'E24' '  -14.34%  '
Set-TextValue 'D25' '13.14'
Set-TextValue 'E25' '  -16.66%  '
Set-TextValue 'D26' '0.999'
Set-TextValue 'E26' '  -0.12%  '
Set-TextValue 'E27' '  -18.08%  '
Set-TextValue 'D28' '2.14'
Set-TextValue 'E28' '  -17.88%  '
Set-TextValue 'D29' '7.89'
Set-TextValue 'E29' '  -14.23%  '
Set-TextValue 'D30' '28.75'
Set-TextValue 'E30' '  -15.31%  '
Set-TextValue 'D31' '2.53'
Set-TextValue 'E31' '  -16.36%  '
Set-TextValue 'D32' '1.11'
Set-TextValue 'E32' '  -14.90%  '
Set-TextValue 'D33' '527.90'
Set-TextValue 'E33' '  -15.03%  '
Set-TextValue 'D34' '6.45'
Set-TextValue 'E34' '  -21.05%  '
Set-TextValue 'B35' 'FirstDigitalUSD'
Set-TextValue 'C35' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.20%  '
Set-TextValue 'B36' 'NEARProtocol'
Set-TextValue 'C36' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D36' '5.58'
Set-TextValue 'E36' '  -18.67%  '
Set-TextValue 'D37' '52.83'
Set-TextValue 'E37' '  -7.34%  '
Set-TextValue 'D38' '0.0848'
Set-TextValue 'E38' '  -15.59%  '
Set-TextValue 'D39' '0.0410'
Set-TextValue 'E39' '  -16.74%  '
Set-TextValue 'D40' '9.00'
Set-TextValue 'E40' '  -17.02%  '
Set-TextValue 'D41' '0.123'
Set-TextValue 'E41' '  -15.10%  '
Set-TextValue 'D42' '2.876.61'
Set-TextValue 'D43' '2.59'
Set-TextValue 'E43' '  -25.99%  '
Set-TextValue 'B44' 'PEPE'
Set-TextValue 'C44' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D44' '0.0₃0578'
Set-TextValue 'E44' '  -21.87%  '
Set-TextValue 'B45' 'TheGraph'
Set-TextValue 'C45' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D45' '0.257'
Set-TextValue 'E45' '  -17.83%  '
Set-TextValue 'D47' '25.56'
Set-TextValue 'E47' '  -20.86%  '
Set-TextValue 'D48' '2.30'
Set-TextValue 'E48' '  -21.61%  '
Set-TextValue 'D49' '2.06'
Set-TextValue 'E49' '  -19.80%  '
Set-TextValue 'D50' '0.112'
Set-TextValue 'E50' '  -14.28%  '
Set-TextValue 'D51' '122.18'
Set-TextValue 'E51' '  -8.77%  '
